$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of monthly data appended after the existing last row (100).
# Column A carries the same date-format cell style as the rows above it,
# so copy that formatting from the last existing row before writing values.

$ws.Range("A100").Copy($ws.Range("A101"))
$ws.Range("A101").Value = 45748
$ws.Range("B101").Value = 0.266329085116067
$ws.Range("C101").Value = 0.130100142453171

$ws.Range("A100").Copy($ws.Range("A102"))
$ws.Range("A102").Value = 45778
$ws.Range("B102").Value = 0.203326505841299
$ws.Range("C102").Value = 0.136818829332596
